# Update the date line and the division-problem answers in the practice
# table. The table is 5 columns wide; only every 4th row (1, 5, 9, 13, 17)
# actually carries text, the rest are spacer rows.

$d = $word.ActiveDocument

# 1. Date heading: "2024-01-18 Thursday" -> "2024-01-19 Friday"
$d.Paragraphs.Item(1).Range.Text = "2024-01-19 Friday"

# 2. Table cell answers, addressed by (row, column) so there is no
#    ambiguity even where two cells end up sharing identical text.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "93÷6=15, 3"
$t.Cell(1,2).Range.Text  = "64÷8=8, 0"
$t.Cell(1,3).Range.Text  = "38÷4=9, 2"
$t.Cell(1,4).Range.Text  = "83÷5=16, 3"
$t.Cell(1,5).Range.Text  = "83÷9=9, 2"

$t.Cell(5,1).Range.Text  = "62÷6=10, 2"
$t.Cell(5,2).Range.Text  = "42÷4=10, 2"
$t.Cell(5,3).Range.Text  = "74÷8=9, 2"
$t.Cell(5,4).Range.Text  = "13÷6=2, 1"
$t.Cell(5,5).Range.Text  = "96÷8=12, 0"

$t.Cell(9,1).Range.Text  = "17÷2=8, 1"
$t.Cell(9,2).Range.Text  = "53÷4=13, 1"
$t.Cell(9,3).Range.Text  = "65÷5=13, 0"
$t.Cell(9,4).Range.Text  = "87÷5=17, 2"
$t.Cell(9,5).Range.Text  = "74÷8=9, 2"

$t.Cell(13,1).Range.Text = "49÷3=16, 1"
$t.Cell(13,2).Range.Text = "19÷9=2, 1"
$t.Cell(13,3).Range.Text = "80÷3=26, 2"
$t.Cell(13,4).Range.Text = "20÷3=6, 2"
$t.Cell(13,5).Range.Text = "59÷2=29, 1"

$t.Cell(17,1).Range.Text = "14÷6=2, 2"
$t.Cell(17,2).Range.Text = "68÷2=34, 0"
$t.Cell(17,3).Range.Text = "71÷5=14, 1"
$t.Cell(17,4).Range.Text = "83÷9=9, 2"
$t.Cell(17,5).Range.Text = "41÷6=6, 5"
